$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.538.55"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +5.73%  "

$ws.Range("D3").Value = "'1.725.32"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +4.76%  "

$ws.Range("E4").Value = "  +0.05%  "

$ws.Range("D5").Value = "'225.90"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.64%  "

$ws.Range("D6").Value = "'0.5346"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.15%  "

$ws.Range("E7").Value = "  -0.01%  "

$ws.Range("D8").Value = "'0.2676"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.15%  "

$ws.Range("D9").Value = "'0.06593"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +4.81%  "

$ws.Range("D10").Value = "'21.67"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +6.74%  "

$ws.Range("D11").Value = "'0.07711"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.61%  "

$ws.Range("D12").Value = "'4.616"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.95%  "

$ws.Range("D13").Value = "'1.727.41"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +5.14%  "

$ws.Range("D14").Value = "'1.962.44"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +4.79%  "

$ws.Range("D15").Value = "'0.5835"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +4.90%  "

$ws.Range("D16").Value = "'0.0₅8286"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.27%  "

$ws.Range("E17").Value = "  +4.64%  "

$ws.Range("D18").Value = "'27.538.09"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +5.90%  "

$ws.Range("D19").Value = "'219.32"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +14.43%  "

$ws.Range("E20").Value = "  -0.01%  "

$ws.Range("D21").Value = "'4.733"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.07%  "

$ws.Range("E22").Value = "  +2.09%  "

$ws.Range("D23").Value = "'6.093"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +3.41%  "

$ws.Range("E24").Value = "  +0.05%  "

$ws.Range("D25").Value = "'145.98"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.28%  "

$ws.Range("D26").Value = "'1.741"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +14.97%  "

$ws.Range("E27").Value = "  +4.79%  "

$ws.Range("D28").Value = "'7.411"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +3.44%  "

$ws.Range("D29").Value = "'16.60"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +5.08%  "

$ws.Range("D30").Value = "'0.05541"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +3.81%  "

$ws.Range("E31").Value = "  +3.13%  "

$ws.Range("E32").Value = "  +3.71%  "

$ws.Range("D33").Value = "'3.451"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +3.62%  "

$ws.Range("D34").Value = "'1.661"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +7.57%  "

$ws.Range("D35").Value = "'2.858"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +3.00%  "

$ws.Range("D36").Value = "'0.9668"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.91%  "

$ws.Range("D37").Value = "'2.425"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.37%  "

$ws.Range("D38").Value = "'0.5988"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +7.31%  "

$ws.Range("E39").Value = "  +5.16%  "

$ws.Range("D40").Value = "'5.904"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.34%  "

$ws.Range("D41").Value = "'0.8576"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +4.16%  "

$ws.Range("D42").Value = "'1.057.08"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.65%  "

$ws.Range("E43").Value = "  +0.00%  "

$ws.Range("D44").Value = "'101.30"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.62%  "

$ws.Range("D45").Value = "'1.869.05"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +4.82%  "

$ws.Range("E46").Value = "  +3.28%  "

$ws.Range("D47").Value = "'59.06"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +3.52%  "

$ws.Range("D48").Value = "'8.226"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +4.25%  "

$ws.Range("D49").Value = "'0.4460"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +3.56%  "

$ws.Range("D50").Value = "'1.003"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.08%  "

$ws.Range("D51").Value = "'0.05242"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.59%  "
